$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the modified cells keep their text (string) representation
# rather than being auto-converted to numbers/percentages by Excel.
$targetCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "E47")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "273.28"
$ws.Range("E2").Value = "1.21%"
$ws.Range("D3").Value = "26.75"
$ws.Range("E3").Value = "0.14%"
$ws.Range("D4").Value = "4.898"
$ws.Range("E4").Value = "3.90%"
$ws.Range("D5").Value = "0.06332"
$ws.Range("E5").Value = "3.72%"
$ws.Range("D6").Value = "6.929"
$ws.Range("E6").Value = "2.81%"
$ws.Range("D7").Value = "3.348"
$ws.Range("E7").Value = "5.52%"
$ws.Range("D8").Value = "1.343"
$ws.Range("E8").Value = "50.76%"
$ws.Range("D9").Value = "0.8832"
$ws.Range("E9").Value = "3.08%"
$ws.Range("D10").Value = "0.1471"
$ws.Range("E10").Value = "2.82%"
$ws.Range("D11").Value = "0.05070"
$ws.Range("E11").Value = "2.43%"
$ws.Range("D12").Value = "0.07374"
$ws.Range("E12").Value = "3.61%"
$ws.Range("D13").Value = "0.03172"
$ws.Range("E13").Value = "-0.39%"
$ws.Range("D14").Value = "0.09030"
$ws.Range("E14").Value = "-0.03%"
$ws.Range("D15").Value = "0.001556"
$ws.Range("E15").Value = "1.01%"
$ws.Range("D16").Value = "0.0006298"
$ws.Range("E16").Value = "3.51%"
$ws.Range("D17").Value = "0.006022"
$ws.Range("E17").Value = "1.30%"
$ws.Range("D18").Value = "3.470"
$ws.Range("E18").Value = "0.22%"
$ws.Range("D19").Value = "2.284"
$ws.Range("E19").Value = "0.88%"
$ws.Range("D20").Value = "0.3143"
$ws.Range("E20").Value = "1.72%"
$ws.Range("D21").Value = "0.1332"
$ws.Range("E21").Value = "2.50%"
$ws.Range("D22").Value = "3.944"
$ws.Range("E22").Value = "2.55%"
$ws.Range("D23").Value = "0.04346"
$ws.Range("E23").Value = "2.44%"
$ws.Range("D24").Value = "0.001178"
$ws.Range("E24").Value = "-0.12%"
$ws.Range("D25").Value = "0.003652"
$ws.Range("E25").Value = "-11.99%"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").Value = "0.03%"
$ws.Range("D27").Value = "0.0001696"
$ws.Range("E27").Value = "0.92%"
$ws.Range("D40").Value = "0.04049"
$ws.Range("E40").Value = "2.60%"
$ws.Range("D41").Value = "0.006621"
$ws.Range("E41").Value = "58.21%"
$ws.Range("D42").Value = "0.1165"
$ws.Range("E42").Value = "4.03%"
$ws.Range("D43").Value = "0.002222"
$ws.Range("E43").Value = "9.07%"
$ws.Range("D44").Value = "0.01260"
$ws.Range("E44").Value = "5.83%"
$ws.Range("D45").Value = "0.00005333"
$ws.Range("E45").Value = "4.12%"
$ws.Range("E46").Value = "148.76%"
$ws.Range("D47").Value = "0.02118"
$ws.Range("E47").Value = "-13.47%"

foreach ($addr in $targetCells) {
    $ws.Range($addr).Style = "Normal"
}
